$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: E2/F2 were stored as text ("code"/"quantity" typed as text in
#     the form); fix them to be real numbers.
$ws.Range("E2").Value = 21312312
$ws.Range("F2").Value = 12312

# --- Row 3: a new stock entry submitted through the app's form. Every
#     field from that form is plain text (even the numeric-looking ones),
#     matching how row 2 was originally captured, and blank fields are
#     stored as empty text cells rather than left completely unset.
$ws.Range("A3").Value = "Lucia"
$ws.Range("B3").Value = "22/09/2022"
$ws.Range("C3").Value = "30/09/2022"

# D3, E3, G3 look numeric but must stay text -> force text via quote
# prefix, then strip the resulting formatting so no explicit style lingers.
$ws.Range("D3").Value = "'12312"
$ws.Range("D3").ClearFormats()

$ws.Range("E3").Value = "'21312"
$ws.Range("E3").ClearFormats()

$ws.Range("G3").Value = "'12313"
$ws.Range("G3").ClearFormats()

# F3, H3, I3 are blank text cells (quote-prefix with nothing after it
# yields an empty-string text value), again with the formatting reset.
$ws.Range("F3").Value = "'"
$ws.Range("F3").ClearFormats()

$ws.Range("H3").Value = "'"
$ws.Range("H3").ClearFormats()

$ws.Range("I3").Value = "'"
$ws.Range("I3").ClearFormats()
